# Apply the "improve bibliography and string handling" edit:
#  - Update L2 note text
#  - Add two new rows of metadata (Black Storm-Petrel, Gentry's Giant-Skipper)
#  - Add hyperlinks on the new ContactEmail cells, matching K2's style
#  - Resize columns B, C, L
#  - Move the active selection to L4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 2: only the Notes cell (L2) text actually changes; C2/F2 keep their
# text ("oceahomo" / "Element Occurrences (EOs)") so nothing else to touch.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 12).Value = "a simple selection of 1 km grid cells containing the breeding EOs"

# ---------------------------------------------------------------------------
# Row 3: Oceanodroma melania / Black Storm-Petrel
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value  = "Oceanodroma melania"
$ws.Cells.Item(3, 2).Value  = "Black Storm-Petrel"
$ws.Cells.Item(3, 3).Value  = "oceamela"
$ws.Cells.Item(3, 4).Value  = 105368
$ws.Cells.Item(3, 5).Value  = "G3"
$ws.Cells.Item(3, 6).Value  = "Element Occurrences (EOs)"
$ws.Cells.Item(3, 7).Value  = "California Natural Diversity Database"
$ws.Cells.Item(3, 8).Value  = "oceamela_AltMap_20190228.tif"
$ws.Cells.Item(3, 9).Value  = "Michelle M. Fink"
$ws.Cells.Item(3, 10).Value = "Colorado Natural Heritage Program"
$ws.Cells.Item(3, 11).Value = "michelle.fink@colostate.edu"
$ws.Cells.Item(3, 12).Value = "a simple selection of 1 km grid cells containing the breeding EOs"
$ws.Cells.Item(3, 13).Value = "Awaiting review"
$ws.Cells.Item(3, 14).Value = "Not final"

# ---------------------------------------------------------------------------
# Row 4: Agathymus gentryi / Gentry's Giant-Skipper
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value  = "Agathymus gentryi"
$ws.Cells.Item(4, 2).Value  = "Gentry's Giant-Skipper"
$ws.Cells.Item(4, 3).Value  = "agatgent"
$ws.Cells.Item(4, 4).Value  = 116354
$ws.Cells.Item(4, 5).Value  = "G3"
$ws.Cells.Item(4, 6).Value  = "Element Occurrences (EOs), Butterflies and Moths of North America (BAMONA), Land Cover, and Elevation"
$ws.Cells.Item(4, 7).Value  = "Arizona Heritage Data Management System"
$ws.Cells.Item(4, 8).Value  = "agatgent_AltMap_20190417.tif"
$ws.Cells.Item(4, 9).Value  = "Michelle M. Fink"
$ws.Cells.Item(4, 10).Value = "Colorado Natural Heritage Program"
$ws.Cells.Item(4, 11).Value = "michelle.fink@colostate.edu"
$ws.Cells.Item(4, 12).Value = "created using EOs plus BAMONA locations to inform land cover (pinyon-juniper and desert scrub) and topography (intermediately to moderately rugged) selections"
$ws.Cells.Item(4, 13).Value = "Awaiting review"
$ws.Cells.Item(4, 14).Value = "Not final"

# ---------------------------------------------------------------------------
# Hyperlinks for the new ContactEmail cells (K3, K4), matching K2's target.
# Adding a hyperlink resets the cell style, so re-apply K2's format after.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("K3"), "mailto:michelle.fink@colostate.edu")
$ws.Range("K2").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("K4"), "mailto:michelle.fink@colostate.edu")
$ws.Range("K2").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column width adjustments (B, C, L)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.3
$ws.Columns.Item(3).ColumnWidth = 9.45
$ws.Columns.Item(12).ColumnWidth = 67.65

# ---------------------------------------------------------------------------
# Move the active selection to L4
# ---------------------------------------------------------------------------
[void]$ws.Range("L4").Select()
